$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.530.55"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "1.639.89"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Cells.Item(4, 4).Value = "'0.9998"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "'308.26"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Cells.Item(7, 4).Value = "'0.3774"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Cells.Item(8, 4).Value = "'53.09"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("E9").Value = "  +1.96%  "
$ws.Cells.Item(10, 4).Value = "'1.280"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +2.58%  "
$ws.Cells.Item(11, 4).Value = "'0.08219"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +3.80%  "
$ws.Cells.Item(14, 4).Value = "'6.684"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E15").Value = "  +3.61%  "
$ws.Cells.Item(16, 4).Value = "'7.491"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "1.640.91"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Cells.Item(18, 4).Value = "'95.17"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Cells.Item(19, 4).Value = "'0.06958"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Cells.Item(20, 4).Value = "'18.45"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Cells.Item(21, 4).Value = "'6.604"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Cells.Item(22, 4).Value = "'0.9979"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "23.535.27"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Cells.Item(25, 4).Value = "'3.112"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +6.33%  "
$ws.Cells.Item(26, 4).Value = "'2.419"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Cells.Item(27, 4).Value = "'21.42"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Cells.Item(28, 4).Value = "'151.32"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Cells.Item(29, 4).Value = "'5.333"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Cells.Item(30, 4).Value = "'136.17"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Cells.Item(31, 4).Value = "'2.426"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'6.911"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("D33").Value = "1.820.82"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Cells.Item(34, 4).Value = "'0.9814"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Cells.Item(35, 4).Value = "'0.02816"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +4.70%  "
$ws.Cells.Item(36, 4).Value = "'10.47"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Cells.Item(37, 4).Value = "'0.07477"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +2.14%  "
$ws.Cells.Item(39, 4).Value = "'0.2552"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +1.96%  "
$ws.Cells.Item(40, 4).Value = "'0.08900"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("E41").Value = "  +3.24%  "
$ws.Cells.Item(42, 4).Value = "'0.7182"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Cells.Item(43, 4).Value = "'12.66"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Cells.Item(44, 4).Value = "'16.38"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  +9.44%  "
$ws.Cells.Item(45, 4).Value = "'0.6634"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("E46").Value = "  +3.83%  "
$ws.Cells.Item(47, 4).Value = "'4.047"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Cells.Item(48, 4).Value = "'0.9990"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Cells.Item(49, 4).Value = "'0.08075"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Cells.Item(50, 4).Value = "'131.05"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Cells.Item(51, 4).Value = "'1.225"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  +0.96%  "
